$d = $word.ActiveDocument

# 1) "Department of Famil" + "y Medicine..." -> merge text (fix accidental run split)
$d.Content.Find.Execute("Department of Famil" + "y Medicine, University of Ottawa, Ottawa, Ontario, Canada", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Department of Family Medicine, University of Ottawa, Ottawa, Ontario, Canada", 2)

# 2) "Civ" + "ic Campus, ASB 2-012" -> merge text (fix accidental run split)
$d.Content.Find.Execute("Civic Campus, ASB 2-012", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Civic Campus, ASB 2-012", 2)

# 3) " The authors d" + "eclare..." -> merge text (fix accidental run split)
$d.Content.Find.Execute("The authors declare that they have no conflict of interest.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "The authors declare that they have no conflict of interest.", 2)

# 4) Word count: Abstract 228 -> 222 ; Text body 2835 -> 3118
$d.Content.Find.Execute("Abstract: 228; Text body: 2835", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Abstract: 222; Text body: 3118", 2)

# 5) The "Pages:" paragraph is removed (its label is reused to replace the old
#    "Figures:" paragraph that immediately follows it), and the References
#    count changes from 30 to 27.
#    Find the "Pages:" paragraph by its exact original text, then delete the
#    paragraph that immediately follows it (the original "Figures: 6" line),
#    and finally relabel the "Pages:" run itself to "Figures:".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "Pages: 6`r") {
        $next = $d.Paragraphs($i + 1)
        if ($next.Range.Text -eq "Figures: 6`r") {
            $next.Range.Delete()
        }
        break
    }
}

$d.Content.Find.Execute("Pages:", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Figures:", 2)

$d.Content.Find.Execute(" 30", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " 27", 2)
